$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to machine-friendly names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case "de"/"de la"/"de los" in place names ---
$ws.Range("B5").Value  = "Comitán De Domínguez"
$ws.Range("B9").Value  = "Mazapa De Madero"
$ws.Range("A14").Value = "Ciudad De México"
$ws.Range("A16").Value = "Estado De México"
$ws.Range("B16").Value = "Almoloya De Alquisiras"
$ws.Range("B18").Value = "San Luis De La Paz"
$ws.Range("B21").Value = "Atoyac De Álvarez"
$ws.Range("B22").Value = "Chilpancingo De Los Bravo"
$ws.Range("B26").Value = "Técpan De Galeana"
$ws.Range("B29").Value = "San Cristóbal De La Barranca"
$ws.Range("A32").Value = "Michoacán De Ocampo"
$ws.Range("B43").Value = "Oaxaca De Juárez"
$ws.Range("B52").Value = "Amealco De Bonfil"
$ws.Range("A57").Value = "Veracruz De Ignacio De La Llave"

# --- Grand total label ---
$ws.Range("A66").Value = "Total"

# --- Remove trailing metadata/footer rows (68-72) ---
$ws.Rows("68:72").Delete()
